# Correct inversion of ZP D (zp_Deep) and ZP S (zp_Shallow) blocks.
# Rows 758-793 hold the "zp_Deep" data (D1/D2/D3 reps) and rows 794-829
# hold the "zp_Shallow" data (S1/S2/S3 reps). The Abs_cov (B), Rel_cov (C)
# and sd_rel_cov (J) columns for these two blocks were swapped in the
# source data; this script swaps them back so each block holds the
# correct values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$deepStart = 758
$deepEnd = 793
$offset = 36   # shallow row = deep row + 36

for ($r = $deepStart; $r -le $deepEnd; $r++) {
    $sRow = $r + $offset

    # Column B - Abs_cov
    $bDeep = $ws.Cells.Item($r, 2).Value()
    $bShallow = $ws.Cells.Item($sRow, 2).Value()
    $ws.Cells.Item($r, 2).Value = $bShallow
    $ws.Cells.Item($sRow, 2).Value = $bDeep

    # Column C - Rel_cov
    $cDeep = $ws.Cells.Item($r, 3).Value()
    $cShallow = $ws.Cells.Item($sRow, 3).Value()
    $ws.Cells.Item($r, 3).Value = $cShallow
    $ws.Cells.Item($sRow, 3).Value = $cDeep

    # Column J - sd_rel_cov
    $jDeep = $ws.Cells.Item($r, 10).Value()
    $jShallow = $ws.Cells.Item($sRow, 10).Value()
    $ws.Cells.Item($r, 10).Value = $jShallow
    $ws.Cells.Item($sRow, 10).Value = $jDeep
}
